$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Lipton Chicken Noodle Soup " ingredient to "Lipton Envelope" (row 10)
$ws.Range("A10").Value = "Lipton Envelope"

# 2. Insert a new row for "Mushrooms" just below "Basil" (row 30), above "Spaghetti" (row 31)
$ws.Rows("31").Insert()
$ws.Range("A31").Value = "Mushrooms"
$ws.Range("B31").Value = "Vegetables"

# 3. Insert a new row for "Tomato" just above the existing "Tomatoes" row (row 29)
$ws.Rows("29").Insert()
$ws.Range("A29").Value = "Tomato"
$ws.Range("B29").Value = "Vegetables"

# 4. Update the view state (scroll position / selection) to match the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("G25").Select()
